# Actualiza base de datos EC: reordena los periodos de mora (columna E,
# filas 16-44) de orden descendente a orden ascendente.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periodos = @(
  "1702","1706","1801","1802","1803","1804","1805","1806","1807","1808",
  "1809","1810","1811","1812","1901","1902","1903","1904","1905","1906",
  "1907","1908","1909","1910","1911","1912","2001","2002","2003"
)

$startRow = 16
for ($i = 0; $i -lt $periodos.Length; $i++) {
  $row = $startRow + $i
  $ws.Range("E$row").Value = $periodos[$i]
}
